# Add a new "Spain" worksheet, cloned from the "Italy" worksheet, with
# Spain-specific market data (per commit: "Added Test data for Spain Zettler Market").

$wb = $excel.ActiveWorkbook

$italy = $wb.Worksheets.Item("Italy")

# Copy Italy to create the new Spain sheet right after it.
$italy.Copy($null, $italy)

$spain = $wb.Worksheets.Item($italy.Index + 1)
$spain.Name = "Spain"

# Update the market-specific values on the new sheet.
# (Order matters for shared-string allocation: the NGC code string is
# appended before the market-name string, matching the authored workbook.)
$spain.Range("B4").Value = "NGC-3103/T2041"
$spain.Range("B2").Value = "Spain Market"

# Selection / active-cell bookkeeping to match the target state.
$italy.Range("A1:D11").Select()
$spain.Activate()
$spain.Range("D8").Select()
